# "Added scroll wheel for drive mode"
#
# Sheet "Translation" (sheet2.xml) lists one row per translatable text id.
# Row 5 ("SingleUseId2": Large / Center / LTR / "<value>") is removed, and
# row 18 ("SingleUseId15": Large / Left / LTR / "STANDARD") is removed
# (rows below each deletion shift up). Two new rows are appended at the
# bottom of the table with fresh auto-generated ids that carry the same
# content as the two removed rows - these back the new drive-mode scroll
# wheel text entries:
#   - SingleUseId18: Large / Center / LTR / "<value>"
#   - SingleUseId19: Large / Left   / LTR / "STANDARD"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Remove the old row 5 entry ("SingleUseId2").
$ws.Rows.Item(5).Delete()

# After the first delete, the old row 18 entry ("SingleUseId15") is now at
# row 17 - remove it too.
$ws.Rows.Item(17).Delete()

# Append the two replacement rows at the (new) end of the table, 19 and 20.
$ws.Range("B19").Value = "SingleUseId18"
$ws.Range("C19").Value = "Large"
$ws.Range("D19").Value = "Center"
$ws.Range("E19").Value = "LTR"
$ws.Range("F19").Value = "<value>"

$ws.Range("B20").Value = "SingleUseId19"
$ws.Range("C20").Value = "Large"
$ws.Range("D20").Value = "Left"
$ws.Range("E20").Value = "LTR"
$ws.Range("F20").Value = "STANDARD"
